# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (E2) and
# "Correspond Handback DateTime" (H2) values on the zh-cn and de-de
# report sheets to reflect the latest handback run timestamps.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-14 08:48:04"
$wsZhCn.Range("H2").Value = "2016-03-14 08:48:21"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-14 08:48:07"
$wsDeDe.Range("H2").Value = "2016-03-14 08:48:26"
